$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = ".victorinsurance.nl"
$ws.Range("C2").Value = "2024-09-05 22:12:04"

$ws.Range("B3").Value = ".victorinsurance.nl"
$ws.Range("C3").Value = "2024-09-05 22:12:04"
